$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Submit orders" -> add row 91 ---
$ws1 = $wb.Worksheets.Item("Submit orders")
$ws1.Cells.Item(91, 1).Value = "10.28.2022 21:49 (Kyiv+Israel) 18:49 (UTC) 03:49 (Japan) 00:19 (India)"
$ws1.Cells.Item(91, 2).Value = 0.835
$ws1.Cells.Item(91, 3).Value = -0.108
$ws1.Cells.Item(91, 4).Value = "***"
$ws1.Cells.Item(91, 5).Value = "***"

# --- Sheet 2: "Submit internet survey" -> add rows 84-89 ---
$ws2 = $wb.Worksheets.Item("Submit internet survey")

$ws2.Cells.Item(84, 1).Value = "10.28.2022 13:00 (Kyiv+Israel) 10:00 (UTC) 19:00 (Japan) 15:30 (India)"
$ws2.Cells.Item(84, 2).Value = "***"
$ws2.Cells.Item(84, 3).Value = "***"
$ws2.Cells.Item(84, 4).Value = 1.023
$ws2.Cells.Item(84, 5).Value = -0.4299999999999999

$ws2.Cells.Item(85, 1).Value = "10.28.2022 13:03 (Kyiv+Israel) 10:03 (UTC) 19:03 (Japan) 15:33 (India)"
$ws2.Cells.Item(85, 2).Value = 1.104
$ws2.Cells.Item(85, 3).Value = -0.4730000000000001
$ws2.Cells.Item(85, 4).Value = "***"
$ws2.Cells.Item(85, 5).Value = "***"

$ws2.Cells.Item(86, 1).Value = "10.28.2022 13:06 (Kyiv+Israel) 10:06 (UTC) 19:06 (Japan) 15:36 (India)"
$ws2.Cells.Item(86, 2).Value = 0.643
$ws2.Cells.Item(86, 3).Value = -0.01200000000000001
$ws2.Cells.Item(86, 4).Value = "***"
$ws2.Cells.Item(86, 5).Value = "***"

$ws2.Cells.Item(87, 1).Value = "10.28.2022 13:21 (Kyiv+Israel) 10:21 (UTC) 19:21 (Japan) 15:51 (India)"
$ws2.Cells.Item(87, 2).Value = 0.979
$ws2.Cells.Item(87, 3).Value = -0.348
$ws2.Cells.Item(87, 4).Value = "***"
$ws2.Cells.Item(87, 5).Value = "***"

$ws2.Cells.Item(88, 1).Value = "10.28.2022 21:06 (Kyiv+Israel) 18:06 (UTC) 03:06 (Japan) 23:36 (India)"
$ws2.Cells.Item(88, 2).Value = 0.716
$ws2.Cells.Item(88, 3).Value = -0.08499999999999996
$ws2.Cells.Item(88, 4).Value = "***"
$ws2.Cells.Item(88, 5).Value = "***"

$ws2.Cells.Item(89, 1).Value = "10.28.2022 21:52 (Kyiv+Israel) 18:52 (UTC) 03:52 (Japan) 00:22 (India)"
$ws2.Cells.Item(89, 2).Value = 0.58
$ws2.Cells.Item(89, 3).Value = 0.05100000000000005
$ws2.Cells.Item(89, 4).Value = "***"
$ws2.Cells.Item(89, 5).Value = "***"

# --- Sheet 3: "Submit a phone survey" -> add row 82 ---
$ws3 = $wb.Worksheets.Item("Submit a phone survey")
$ws3.Cells.Item(82, 1).Value = "10.28.2022 21:54 (Kyiv+Israel) 18:54 (UTC) 03:54 (Japan) 00:24 (India)"
$ws3.Cells.Item(82, 2).Value = 1.398
$ws3.Cells.Item(82, 3).Value = -0.2939999999999998
$ws3.Cells.Item(82, 4).Value = "***"
$ws3.Cells.Item(82, 5).Value = "***"

# --- Sheet 4: "Checkertificate" -> add row 95 ---
$ws4 = $wb.Worksheets.Item("Checkertificate")
$ws4.Cells.Item(95, 1).Value = "10.28.2022 22:02 (Kyiv+Israel) 19:02 (UTC) 04:02 (Japan) 00:32 (India)"
$ws4.Cells.Item(95, 2).Value = 0.733
$ws4.Cells.Item(95, 3).Value = -0.06799999999999995
$ws4.Cells.Item(95, 4).Value = "***"
$ws4.Cells.Item(95, 5).Value = "***"
